# Add a new "case" worksheet at the end of the workbook (becomes the active/selected tab),
# populate it with the test-case table, and remove the previous tab selection from "forms".

$wb = $excel.ActiveWorkbook

# Insert the new worksheet right after the current last sheet ("forms") so it lands at the end.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "case"

# Header row
$newSheet.Cells.Item(1,1).Value = "title"
$newSheet.Cells.Item(1,2).Value = "state"
$newSheet.Cells.Item(1,3).Value = "identifier"
$newSheet.Cells.Item(1,4).Value = "tags"
$newSheet.Cells.Item(1,5).Value = "description"
$newSheet.Cells.Item(1,6).Value = "type"
$newSheet.Cells.Item(1,7).Value = "priority"

# Row 2 - fill columns A through E first
$newSheet.Cells.Item(2,1).Value = "Test Case title - 1"
$newSheet.Cells.Item(2,2).Value = "Enquiring"
$newSheet.Cells.Item(2,3).Value = "Test identifier - 1"
$newSheet.Cells.Item(2,4).Value = "Test tags - 1"
$newSheet.Cells.Item(2,5).Value = "Test description - 1"

# Row 3 - state (B) entered before title (A) and the rest
$newSheet.Cells.Item(3,2).Value = "Reviewing"
$newSheet.Cells.Item(3,1).Value = "Test Case title - 2"
$newSheet.Cells.Item(3,3).Value = "Test identifier - 2"
$newSheet.Cells.Item(3,4).Value = "Test tags - 2"
$newSheet.Cells.Item(3,5).Value = "Test description - 2"

# Type column (F) filled last for both rows
$newSheet.Cells.Item(2,6).Value = "Complaint"
$newSheet.Cells.Item(3,6).Value = "General Support"

# Priority column (G) reuses existing values
$newSheet.Cells.Item(2,7).Value = "High"
$newSheet.Cells.Item(3,7).Value = "Normal"

# Column widths (best fit, matching the authored layout)
$newSheet.Columns.Item(1).ColumnWidth = 15.307291666666666
$newSheet.Columns.Item(2).ColumnWidth = 8.592447916666666
$newSheet.Columns.Item(3).ColumnWidth = 15.451822916666666
$newSheet.Columns.Item(4).ColumnWidth = 10.451822916666666
$newSheet.Columns.Item(5).ColumnWidth = 17.166666666666668
$newSheet.Columns.Item(6).ColumnWidth = 14.736979166666666
$newSheet.Columns.Item(7).ColumnWidth = 6.736979166666667

# Selection on the new sheet
$newSheet.Range("G5").Select()
